$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": update changed cells per latest BRVM data refresh ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$ws1.Cells.Item(2,4).Value = 2433.32
$ws1.Cells.Item(2,5).Value = 100.02
$ws1.Cells.Item(3,4).Value = 1914.01
$ws1.Cells.Item(3,5).Value = 635.58
$ws1.Cells.Item(4,1).Value = "UNIWAX CI"
$ws1.Cells.Item(4,4).Value = 1780
$ws1.Cells.Item(4,5).Value = 620
$ws1.Cells.Item(5,1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(5,4).Value = 1745
$ws1.Cells.Item(5,5).Value = 555
$ws1.Cells.Item(6,4).Value = 1680
$ws1.Cells.Item(6,5).Value = 560
$ws1.Cells.Item(7,4).Value = 1640
$ws1.Cells.Item(8,4).Value = 1112.86
$ws1.Cells.Item(8,5).Value = 372.3
$ws1.Cells.Item(9,4).Value = 1069.66
$ws1.Cells.Item(9,5).Value = 367.16
$ws1.Cells.Item(10,1).Value = "SAFCA CI"
$ws1.Cells.Item(10,3).Value = 1
$ws1.Cells.Item(10,4).Value = 995
$ws1.Cells.Item(10,5).Value = 995
$ws1.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(11,4).Value = 975.16
$ws1.Cells.Item(11,5).Value = 322.67
$ws1.Cells.Item(12,1).Value = "BRVM - INDUSTRIE"
$ws1.Cells.Item(12,4).Value = 788.8200000000001
$ws1.Cells.Item(12,5).Value = 263.29
$ws1.Cells.Item(13,4).Value = 665
$ws1.Cells.Item(13,5).Value = 665
$ws1.Cells.Item(14,1).Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Cells.Item(14,4).Value = 652.41
$ws1.Cells.Item(14,5).Value = 217.34
$ws1.Cells.Item(15,1).Value = "BRVM-PRINCIPAL"
$ws1.Cells.Item(15,4).Value = 570.67
$ws1.Cells.Item(15,5).Value = 190.36
$ws1.Cells.Item(16,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(16,4).Value = 413.3
$ws1.Cells.Item(16,5).Value = 139.13
$ws1.Cells.Item(17,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(17,4).Value = 389.4
$ws1.Cells.Item(17,5).Value = 130.83
$ws1.Cells.Item(18,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(18,4).Value = 367.67
$ws1.Cells.Item(18,5).Value = 123.28
$ws1.Cells.Item(19,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(19,4).Value = 361.34
$ws1.Cells.Item(19,5).Value = 121.16
$ws1.Cells.Item(20,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(20,4).Value = 331.8
$ws1.Cells.Item(20,5).Value = 110.89
$ws1.Cells.Item(21,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(21,4).Value = 318.46
$ws1.Cells.Item(21,5).Value = 106.39
$ws1.Cells.Item(22,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(22,2).Value = 0
$ws1.Cells.Item(22,3).Value = 3
$ws1.Cells.Item(22,4).Value = 278.57
$ws1.Cells.Item(22,5).Value = 93.19
$ws1.Cells.Item(22,7).Value = "➖ Neutre"
$ws1.Cells.Item(23,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(23,2).Value = 2
$ws1.Cells.Item(23,4).Value = 10.4
$ws1.Cells.Item(23,5).Value = 7.07
$ws1.Cells.Item(24,1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(24,2).Value = 2
$ws1.Cells.Item(24,3).Value = 1
$ws1.Cells.Item(24,4).Value = 10.21
$ws1.Cells.Item(24,5).Value = 6.25
$ws1.Cells.Item(24,7).Value = "👀 À surveiller"
$ws1.Cells.Item(25,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 6.9
$ws1.Cells.Item(25,5).Value = 6.9
$ws1.Cells.Item(25,7).Value = "➖ Neutre"
$ws1.Cells.Item(26,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(26,2).Value = 2
$ws1.Cells.Item(26,3).Value = 1
$ws1.Cells.Item(26,4).Value = 6.62
$ws1.Cells.Item(26,5).Value = 6.25
$ws1.Cells.Item(26,7).Value = "👀 À surveiller"
$ws1.Cells.Item(27,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(27,4).Value = 5.73
$ws1.Cells.Item(27,5).Value = 5.73
$ws1.Cells.Item(28,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(28,4).Value = 4.03
$ws1.Cells.Item(28,5).Value = 4.03
$ws1.Cells.Item(29,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(29,4).Value = 3.81
$ws1.Cells.Item(29,5).Value = 3.81
$ws1.Cells.Item(30,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(30,4).Value = 3.54
$ws1.Cells.Item(30,5).Value = 3.54
$ws1.Cells.Item(31,1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(31,4).Value = 3.05
$ws1.Cells.Item(31,5).Value = 3.05
$ws1.Cells.Item(32,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(32,4).Value = 2.93
$ws1.Cells.Item(32,5).Value = 2.93
$ws1.Cells.Item(34,1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(34,4).Value = 0.44
$ws1.Cells.Item(34,5).Value = -2.25
$ws1.Cells.Item(37,1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(37,4).Value = -1.87
$ws1.Cells.Item(37,5).Value = -1.87
$ws1.Cells.Item(38,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(38,4).Value = -1.97
$ws1.Cells.Item(38,5).Value = -1.97
$ws1.Cells.Item(39,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(39,4).Value = -2.02
$ws1.Cells.Item(39,5).Value = -2.02
$ws1.Cells.Item(40,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(40,4).Value = -3.16
$ws1.Cells.Item(40,5).Value = -3.16
$ws1.Cells.Item(41,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(41,4).Value = -3.68
$ws1.Cells.Item(41,5).Value = -3.68
$ws1.Cells.Item(42,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(42,3).Value = 2
$ws1.Cells.Item(42,4).Value = -5.97
$ws1.Cells.Item(42,5).Value = -5.37
$ws1.Cells.Item(43,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(43,3).Value = 3
$ws1.Cells.Item(43,4).Value = -21.97
$ws1.Cells.Item(43,5).Value = -7.44
$ws1.Cells.Item(43,6).Value = "🔴 Vente"
$ws1.Cells.Item(43,7).Value = "⚠️ Risque de décrochage"

# Row 44 (SAFCA CI (SAFC)) no longer present as a standalone trailing row; remove it
$ws1.Rows.Item(44).Delete()

# --- Sheet "Top_YTD": refresh values / reordered labels ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Cells.Item(2,2).Value = 426555.13
$ws2.Cells.Item(3,2).Value = 40094.78
$ws2.Cells.Item(4,1).Value = "UNIWAX CI"
$ws2.Cells.Item(4,2).Value = 33192.8
$ws2.Cells.Item(5,1).Value = "NEI-CEDA CI"
$ws2.Cells.Item(5,2).Value = 31538.14
$ws2.Cells.Item(6,2).Value = 28623.2
$ws2.Cells.Item(7,2).Value = 26941.63
$ws2.Cells.Item(8,2).Value = 10345.44
$ws2.Cells.Item(9,2).Value = 9412.34
$ws2.Cells.Item(10,2).Value = 7579.26
$ws2.Cells.Item(11,2).Value = 4680.82
